$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 5
